# Apply the StructureDefinition-insight-type.xlsx update:
#  - Metadata sheet: bump Version to 6.0.0, update Date, fill in Publisher,
#    replace the duplicate "Contact" row with "Jurisdiction", and delete the
#    now-redundant extra "Contact" row (sheet shrinks from 21 to 20 rows).
#  - Elements sheet: update the root Extension row's Short/Definition text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 ("Contact" / "No display for ContactDetail") becomes Jurisdiction
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row - remove it entirely, shifting rows 12+ up
$meta.Rows.Item(11).Delete()

# Elements sheet: root Extension element's Short / Definition text
$elements.Range("K2").Value = "Insight Type"
$elements.Range("L2").Value = "Code for the specific insight type."
